# Append 34 new data rows (580-613) to the "data" sheet and refresh the
# describe()-style summary statistics on the "desc_stat" sheet.

$wb = $excel.ActiveWorkbook
$wsData = $wb.Worksheets.Item("data")
$wsStat = $wb.Worksheets.Item("desc_stat")

# New observations (columns: y, x1, x2, x3)
$newData = @(
    @(14.34854534416843, 3.965, 6.03, 0.296),
    @(8.001250000001374, 4.98, 5.89, 0.3181176470588235),
    @(5.912001383529192, 4.29, 5.94, 0.2710999999999421),
    @(6.68, 4.015, 4.23, 0.324),
    @(7.58, 4.435, 4.53, 0.353206030949981),
    @(11.62562499999963, 5.01, 4.11, 0.338),
    @(13.79999999999967, 4.365, 3.42, 0.333),
    @(10.90750000000006, 4.665, 4.23, 0.3018664137845254),
    @(4.52, 4.005, 4.56, 0.2369987287018682),
    @(5.978823529412008, 4.665, 6.8, 0.3371642162301796),
    @(4.44, 4.895, 6.53, 0.4232575859657763),
    @(7.62, 4.865, 7.44, 0.4129768452924647),
    @(10.67126259564734, 4.83, 5.13, 0.3671176470588323),
    @(10.79812499999665, 4.91, 3.69, 0.449),
    @(10.35, 4.78, 4.88, 0.387),
    @(6.17, 4.91, 7.52, 0.368),
    @(4.277868628312075, 5.61, 2.69, 0.362),
    @(13.64687499999914, 2.62, 2.8, 0.375),
    @(19.79437500000008, 3.09, 3.14, 0.325),
    @(6.617159854848987, 4.585, 5.05, 0.321),
    @(5.11, 5.115, 4.1, 0.3278128355862567),
    @(4.29, 4.635, 4.31, 0.313),
    @(6.11, 5.529999999999999, 5.31, 0.324748561669166),
    @(6.572336983837038, 4.695, 4.01, 0.277),
    @(7.38, 5.205, 6.1, 0.37),
    @(12.26937499999959, 11.4, 2.95, 0.289687824629983),
    @(8.699999999999999, 3.92, 4.33, 0.3743933793495576),
    @(11.43, 5.085, 3.15, 0.371),
    @(14.95187500000067, 4.76, 3.06, 0.4124117647058823),
    @(11.63680753643135, 5.675000000000001, 2.27, 0.361470588235294),
    @(16.25602045285562, 4.65, 6.31, 0.2826468362297734),
    @(6.51249994650501, 21.56, 4.31, 0.3110614830199199),
    @(20.15687499999154, 4.79, 5.38, 0.37),
    @(26.03, 9.529999999999999, 2.67, 0.3958235294117559)
)

$startRow = 580
for ($i = 0; $i -lt $newData.Length; $i++) {
    $row = $startRow + $i
    $rowData = $newData[$i]
    $wsData.Cells.Item($row, 1).Value = $rowData[0]
    $wsData.Cells.Item($row, 2).Value = $rowData[1]
    $wsData.Cells.Item($row, 3).Value = $rowData[2]
    $wsData.Cells.Item($row, 4).Value = $rowData[3]
}

# Refresh the descriptive statistics table (row 2 = count, 3 = mean, 4 = std,
# 5 = min, 6 = 25%, 7 = 50%, 8 = 75%, 9 = max) for columns B:E (y, x1, x2, x3)
$wsStat.Cells.Item(2, 2).Value = 612
$wsStat.Cells.Item(2, 3).Value = 612
$wsStat.Cells.Item(2, 4).Value = 612
$wsStat.Cells.Item(2, 5).Value = 612

$wsStat.Cells.Item(3, 2).Value = 10.77136470793388
$wsStat.Cells.Item(3, 3).Value = 4.553186274509804
$wsStat.Cells.Item(3, 4).Value = 5.030539215686274
$wsStat.Cells.Item(3, 5).Value = 0.3525667678396732

$wsStat.Cells.Item(4, 2).Value = 5.609864831578292
$wsStat.Cells.Item(4, 3).Value = 4.048237262985184
$wsStat.Cells.Item(4, 4).Value = 1.788318488084889
$wsStat.Cells.Item(4, 5).Value = 0.04023879923701083

# row 5 (min) is unchanged

$wsStat.Cells.Item(6, 2).Value = 6.3775
$wsStat.Cells.Item(6, 3).Value = 3.53
$wsStat.Cells.Item(6, 4).Value = 3.6975
$wsStat.Cells.Item(6, 5).Value = 0.324

$wsStat.Cells.Item(7, 2).Value = 9.035
$wsStat.Cells.Item(7, 3).Value = 5.03
$wsStat.Cells.Item(7, 4).Value = 4.655
# row 7 column E (50%) is unchanged (0.351)

$wsStat.Cells.Item(8, 2).Value = 13.80499999999975
$wsStat.Cells.Item(8, 3).Value = 5.765000000000001
$wsStat.Cells.Item(8, 4).Value = 6.1625
$wsStat.Cells.Item(8, 5).Value = 0.381

# row 9 (max) is unchanged
